$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "mała kaczka"
$ws.Range("B4").Value = "duża kaczka"
$ws.Range("C4").Value = "mała kaczka"
